$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New match row appended at the bottom of the table (row 99).
# Copy the formatting of the previous data row (98) down onto the new
# row first so the index column keeps its bold/border style (s="1")
# and the match-date column keeps its custom date-time format (s="2"),
# then overwrite the values.
$ws.Range("A98:V98").Copy()
$ws.Range("A99:V99").PasteSpecial(-4122)

$ws.Range("A99").Value = 98
$ws.Range("B99").Value = "serbia"
$ws.Range("C99").Value = "super-liga"
$ws.Range("D99").Value = "2023-2024"
$ws.Range("E99").Value = 45236.625
$ws.Range("F99").Value = "Novi Pazar"
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = "Partizan"
$ws.Range("I99").Value = 1
$ws.Range("J99").Value = 5.12
$ws.Range("K99").Value = "02/11/2023 08:13"
$ws.Range("L99").Value = 7.33
$ws.Range("M99").Value = "06/11/2023 14:59"
$ws.Range("N99").Value = 3.97
$ws.Range("O99").Value = "02/11/2023 08:13"
$ws.Range("P99").Value = 3.64
$ws.Range("Q99").Value = "06/11/2023 14:52"
$ws.Range("R99").Value = 1.5
$ws.Range("S99").Value = "02/11/2023 08:13"
$ws.Range("T99").Value = 1.52
$ws.Range("U99").Value = "06/11/2023 14:51"
$ws.Range("V99").Value = "https://www.betexplorer.com/football/serbia/super-liga/novi-pazar-partizan/tfzT3o46/"
